# fix: rubber duck and data
# Adds 12 new "duck" rows (rows 10-21) to Sheet1, matching the
# Image_Animal data export used elsewhere in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B (File Name) -- row order 10..21 --------------------------
$ws.Range("B10").Value = "Rita_duck.jpeg"
$ws.Range("B11").Value = "Ivan_duck.jpeg"
$ws.Range("B12").Value = "DanielaS_duck.jpeg"
$ws.Range("B13").Value = "Balduin_duck.jpeg"
$ws.Range("B14").Value = "Christian_duck.jpeg"
$ws.Range("B15").Value = "JulienS_duck.jpeg"
$ws.Range("B16").Value = "Johannes_duck.jpeg"
$ws.Range("B17").Value = "Florian_duck.jpeg"
$ws.Range("B18").Value = "Noemi_duck.jpeg"
$ws.Range("B19").Value = "Nora_duck.jpeg"
$ws.Range("B20").Value = "Gregor_duck.jpeg"
$ws.Range("B21").Value = "RDU_Ginny.jpeg"

# --- Column A (ID) -- row order 10..21 ----------------------------------
$ws.Range("A10").Value = "IA_0009"
$ws.Range("A11").Value = "IA_0010"
$ws.Range("A12").Value = "IA_0011"
$ws.Range("A13").Value = "IA_0012"
$ws.Range("A14").Value = "IA_0013"
$ws.Range("A15").Value = "IA_0014"
$ws.Range("A16").Value = "IA_0015"
$ws.Range("A17").Value = "IA_0016"
$ws.Range("A18").Value = "IA_0017"
$ws.Range("A19").Value = "IA_0018"
$ws.Range("A20").Value = "IA_0019"
$ws.Range("A21").Value = "IA_0020"

# --- Column D (Copyright) -- single value for all new rows -------------
$ws.Range("D10:D21").Value = "Daniela Subotic"

# --- Column E (License List) -- reused existing value -------------------
$ws.Range("E10:E21").Value = "CC BY 4.0"

# --- Column F (Image Directory) -- reused existing value -----------------
$ws.Range("F10:F21").Value = "~/Documents/daschland-scripts/data/Multimedia_Data/Image_Animal/"

# --- Column G (Label) -- row order 10..21 --------------------------------
$ws.Range("G10").Value = "Rita's duck"
$ws.Range("G11").Value = "Ivan's duck"
$ws.Range("G12").Value = "DanielaS's duck"
$ws.Range("G13").Value = "Balduin's duck"
$ws.Range("G14").Value = "Christian's duck"
$ws.Range("G15").Value = "JulienS's duck"
$ws.Range("G16").Value = "Johannes's duck"
$ws.Range("G17").Value = "Florian's duck"
$ws.Range("G18").Value = "Noémi's duck"
$ws.Range("G19").Value = "Nora's duck"
$ws.Range("G20").Value = "Gregor's duck"
$ws.Range("G21").Value = "Ginny in RDU Team meeting"

# --- Column H (Part of Animal Character ID) -- specific entry order ------
$ws.Range("H21").Value = "A_001"
$ws.Range("H10").Value = "A_010"
$ws.Range("H11").Value = "A_011"
$ws.Range("H12").Value = "A_012"
$ws.Range("H13").Value = "A_013"
$ws.Range("H14").Value = "A_014"
$ws.Range("H15").Value = "A_015"
$ws.Range("H16").Value = "A_016"
$ws.Range("H17").Value = "A_017"
$ws.Range("H19").Value = "A_018"
$ws.Range("H20").Value = "A_019"
$ws.Range("H18").Value = "A_009"

# H19/H20 carry a distinct (duplicated) font entry in the source file --
# reapplying the font name nudges the engine into minting a fresh font
# record instead of reusing the shared column style.
$ws.Range("H19").Font.Name = "Aptos Narrow"
$ws.Range("H20").Font.Name = "Aptos Narrow"

# --- Column I (Seqnum) -- numeric literals -------------------------------
$ws.Range("I10").Value = 7
$ws.Range("I11").Value = 8
$ws.Range("I12").Value = 9
$ws.Range("I13").Value = 10
$ws.Range("I14").Value = 11
$ws.Range("I15").Value = 12
$ws.Range("I16").Value = 13
$ws.Range("I17").Value = 14
$ws.Range("I18").Value = 15
$ws.Range("I19").Value = 16
$ws.Range("I20").Value = 17
$ws.Range("I21").Value = 6

# --- Column widths: widen G to fit the longest new label -----------------
$ws.Columns.Item(7).ColumnWidth = 26.83

# --- Selection / view state ----------------------------------------------
$ws.Range("E21").Select()
